# Support TexturePacker. Use UIAtlas to render image from imageset.
# Adds two new columns (K: simple_atlas, L: simple_tex) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Header rows (1-3) ----
$ws.Range("K1").Value = "simple_atlas"
$ws.Range("L1").Value = "simple_tex"

$ws.Range("K2").Value = "string"
$ws.Range("L2").Value = "string"

$ws.Range("K3").Value = "缩略图集"
$ws.Range("L3").Value = "缩略图"

# ---- Data rows (4-20): K = constant "CardSimple", L = same tex id as column H ----
$texIds = @{
    4  = "zhaoyun"
    5  = "lvbu"
    6  = "diaochan"
    7  = "guanyu"
    8  = "zhangfei"
    9  = "simayi"
    10 = "caocao"
    11 = "xiahoudun"
    12 = "yueying"
    13 = "zhugeliang"
    14 = "ganning"
    15 = "shangxiang"
    16 = "huanggai"
    17 = "dianwei"
    18 = "zhanghe"
    19 = "zhoutai"
    20 = "zhouyu"
}

foreach ($row in 4..20) {
    $ws.Range("K$row").Value = "CardSimple"
    $ws.Range("L$row").Value = $texIds[$row]
}

# ---- Column width for the new column K/L boundary (column 11) ----
$ws.Columns.Item(11).ColumnWidth = 19.7

# ---- Update the view selection to match the authored state ----
$ws.Range("N9").Select()
